$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "WALI KELAS" column (F) entirely
$ws.Columns.Item(6).Delete()

# Delete the "NO TELEPON *" column (now shifted to I after the previous delete)
$ws.Columns.Item(9).Delete()

# Update header for column E: "KELAS *" -> "NO TELP *"
$ws.Range("E1").Value = "NO TELP *"

# Update row 2 data
$ws.Range("C2").Value = "dwako"
$ws.Range("E2").Value = 8192391283

# Match the active cell selection shown in the saved file
$ws.Range("I2").Select() | Out-Null
